$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the Password values for rows 2 and 3
$ws.Range("B2").Value = "Test@12345"
$ws.Range("B3").Value = "Test@123"

# Update the active selection shown in the sheet view
$ws.Range("A7").Select()

# Update workbook window position
$excel.ActiveWindow.Left = 1440
$excel.ActiveWindow.Top = 2352
